$wb = $excel.ActiveWorkbook

# 1. Rename the third sheet ("Include from Coverage SelfPay" -> "Include ValueSets")
$wsInclude = $wb.Worksheets.Item("Include from Coverage SelfPay")
$wsInclude.Name = "Include ValueSets"

# 2. Update metadata on the "Metadata" sheet
$wsMeta = $wb.Worksheets.Item("Metadata")

# Date value (row 8, column B)
$wsMeta.Range("B8").Value = "2022-02-22T23:30:09+07:00"

# Description value (row 12, column B) - append the HL7 suffix
$wsMeta.Range("B12").Value = "รหัสประเภทสิทธิ์การรักษาพยาบาลภาครัฐของไทย รวมกับรหัสที่ HL7 กำหนด"

# 3. Replace the content of the "Include ValueSets" sheet:
#    it now only has two rows: a header "ValueSet URL" and the value URL.
#    Remove rows 3 and 4 (and the old System URI / empty rows / column B content).
$wsInclude.Rows.Item(3).Delete()
$wsInclude.Rows.Item(3).Delete()

$wsInclude.Range("A1").Value = "ValueSet URL"
$wsInclude.Range("A2").Value = "http://hl7.org/fhir/ValueSet/coverage-type"
